$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1277066666666667
$ws.Range("H2").Value = 0.38312
$ws.Range("I2").Value = 0.1657851361976475
$ws.Range("J2").Value = 0.1974405699349423
$ws.Range("M2").Value = 0.1798956666666667
$ws.Range("N2").Value = 0.539687
$ws.Range("O2").Value = 0.01070918551864568
$ws.Range("P2").Value = 0.01088280728907136
$ws.Range("Q2").Value = 0.02297387593777778
$ws.Range("R2").Value = 0.20676488344
$ws.Range("S2").Value = 0.001775423779774548
$ws.Range("T2").Value = 0.002148707673646395

# Row 3
$ws.Range("G3").Value = 0.1277066666666667
$ws.Range("H3").Value = 0.38312
$ws.Range("I3").Value = 0.1657851361976475
$ws.Range("J3").Value = 0.1974405699349423
$ws.Range("N3").Value = 47.39813
$ws.Range("O3").Value = 0.9405365839956962
$ws.Range("P3").Value = 0.9557849543390003
$ws.Range("Q3").Value = 2.017685729511111
$ws.Range("R3").Value = 18.1591715656
$ws.Range("S3").Value = 0.1559269856765966
$ws.Range("T3").Value = 0.1887107261199351

# Row 4
$ws.Range("G4").Value = 0.1277066666666667
$ws.Range("H4").Value = 0.38312
$ws.Range("I4").Value = 0.1657851361976475
$ws.Range("J4").Value = 0.1974405699349423
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01499966666666667
$ws.Range("N4").Value = 0.044999
$ws.Range("O4").Value = 0.0008929298633347419
$ws.Range("P4").Value = 0.0009074064137192897
$ws.Range("Q4").Value = 0.001915557431111111
$ws.Range("R4").Value = 0.01724001688
$ws.Range("S4").Value = 0.0001480344990078969
$ws.Range("T4").Value = 0.0001791588394873586

# Row 5
$ws.Range("G5").Value = 0.1277066666666667
$ws.Range("H5").Value = 0.38312
$ws.Range("I5").Value = 0.1657851361976475
$ws.Range("J5").Value = 0.1974405699349423
$ws.Range("M5").Value = 0.8039865
$ws.Range("N5").Value = 1.607973
$ws.Range("O5").Value = 0.04786130062232345
$ws.Range("P5").Value = 0.03242483195820901
$ws.Range("Q5").Value = 0.10267443596
$ws.Range("R5").Value = 0.6160466157600001
$ws.Range("S5").Value = 0.007934692242268442
$ws.Range("T5").Value = 0.00640197730187352

# Row 6
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.370511
$ws.Range("H6").Value = 0.741022
$ws.Range("I6").Value = 0.4809867660085082
$ws.Range("J6").Value = 0.3818850647690823
$ws.Range("M6").Value = 0.1798956666666667
$ws.Range("N6").Value = 0.539687
$ws.Range("O6").Value = 0.01070918551864568
$ws.Range("P6").Value = 0.01088280728907136
$ws.Range("Q6").Value = 0.06665332335233333
$ws.Range("R6").Value = 0.399919940114
$ws.Range("S6").Value = 0.005150976509198533
$ws.Range("T6").Value = 0.004155981566456459

# Row 7
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.370511
$ws.Range("H7").Value = 0.741022
$ws.Range("I7").Value = 0.4809867660085082
$ws.Range("J7").Value = 0.3818850647690823
$ws.Range("N7").Value = 47.39813
$ws.Range("O7").Value = 0.9405365839956962
$ws.Range("P7").Value = 0.9557849543390003
$ws.Range("Q7").Value = 5.853842848143334
$ws.Range("R7").Value = 35.12305708886
$ws.Range("S7").Value = 0.4523856498487795
$ws.Range("T7").Value = 0.3649999991930636

# Row 8
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.370511
$ws.Range("H8").Value = 0.741022
$ws.Range("I8").Value = 0.4809867660085082
$ws.Range("J8").Value = 0.3818850647690823
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.01499966666666667
$ws.Range("N8").Value = 0.044999
$ws.Range("O8").Value = 0.0008929298633347419
$ws.Range("P8").Value = 0.0009074064137192897
$ws.Range("Q8").Value = 0.005557541496333333
$ws.Range("R8").Value = 0.03334524897799999
$ws.Range("S8").Value = 0.0004294874472377967
$ws.Range("T8").Value = 0.0003465249570750717

# Row 9
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.370511
$ws.Range("H9").Value = 0.741022
$ws.Range("I9").Value = 0.4809867660085082
$ws.Range("J9").Value = 0.3818850647690823
$ws.Range("M9").Value = 0.8039865
$ws.Range("N9").Value = 1.607973
$ws.Range("O9").Value = 0.04786130062232345
$ws.Range("P9").Value = 0.03242483195820901
$ws.Range("Q9").Value = 0.2978858421015
$ws.Range("R9").Value = 1.191543368406
$ws.Range("S9").Value = 0.02302065220329235
$ws.Range("T9").Value = 0.01238255905248726

# Row 10
$ws.Range("G10").Value = 0.2720966666666667
$ws.Range("H10").Value = 0.81629
$ws.Range("I10").Value = 0.3532280977938443
$ws.Range("J10").Value = 0.4206743652959753
$ws.Range("M10").Value = 0.1798956666666667
$ws.Range("N10").Value = 0.539687
$ws.Range("O10").Value = 0.01070918551864568
$ws.Range("P10").Value = 0.01088280728907136
$ws.Range("Q10").Value = 0.04894901124777778
$ws.Range("R10").Value = 0.44054110123
$ws.Range("S10").Value = 0.003782785229672597
$ws.Range("T10").Value = 0.00457811804896851

# Row 11
$ws.Range("G11").Value = 0.2720966666666667
$ws.Range("H11").Value = 0.81629
$ws.Range("I11").Value = 0.3532280977938443
$ws.Range("J11").Value = 0.4206743652959753
$ws.Range("N11").Value = 47.39813
$ws.Range("O11").Value = 0.9405365839956962
$ws.Range("P11").Value = 0.9557849543390003
$ws.Range("Q11").Value = 4.298957726411111
$ws.Range("R11").Value = 38.6906195377
$ws.Range("S11").Value = 0.3322239484703201
$ws.Range("T11").Value = 0.4020742290260018

# Row 12
$ws.Range("G12").Value = 0.2720966666666667
$ws.Range("H12").Value = 0.81629
$ws.Range("I12").Value = 0.3532280977938443
$ws.Range("J12").Value = 0.4206743652959753
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.01499966666666667
$ws.Range("N12").Value = 0.044999
$ws.Range("O12").Value = 0.0008929298633347419
$ws.Range("P12").Value = 0.0009074064137192897
$ws.Range("Q12").Value = 0.004081359301111111
$ws.Range("R12").Value = 0.03673223371
$ws.Range("S12").Value = 0.0003154079170890483
$ws.Range("T12").Value = 0.0003817226171568594

# Row 13
$ws.Range("G13").Value = 0.2720966666666667
$ws.Range("H13").Value = 0.81629
$ws.Range("I13").Value = 0.3532280977938443
$ws.Range("J13").Value = 0.4206743652959753
$ws.Range("M13").Value = 0.8039865
$ws.Range("N13").Value = 1.607973
$ws.Range("O13").Value = 0.04786130062232345
$ws.Range("P13").Value = 0.03242483195820901
$ws.Range("Q13").Value = 0.218762046695
$ws.Range("R13").Value = 1.31257228017
$ws.Range("S13").Value = 0.01690595617676265
$ws.Range("T13").Value = 0.01364029560384823

